$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('70328', 'Corn Meal - Coarse', '4', '25.55', '102.20'),
    @('33069', 'Nuts - Almonds (sliced)', '1', '110.03', '110.03'),
    @('33422', 'Nuts - Walnut Halves & Pieces', '1', '125.12', '125.12'),
    @('47170', 'Seeds - Sesame', '8', '67.08', '536.64'),
    @('17617', 'Flour - Whole Wheat (Medium)', '5', '15.89', '79.45'),
    @('18580', 'Flour - Victoria Patent (Org)', '10', '40.02', '400.20'),
    @('17616', 'Flour - Millers Choice', '50', '16.81', '840.50'),
    @('67214-BC', 'Extract - Vanilla', '2', '149.41', '298.82'),
    @('47183', 'Seeds - Sunflower Shelled Raw', '1', '57.08', '57.08'),
    @('28173', 'Sugar - 10X (Powdered)', '1', '41.69', '41.69'),
    @('72020', 'Baking Soda', '1', '37.43', '37.43'),
    @('17297', 'Flour - Cake', '1', '17.10', '17.10'),
    @('70135', 'Condensed Milk', '1', '163.96', '163.96'),
    @('27046', 'Jam - Raspberry (Baking)', '3', '168.31', '504.93'),
    @('39278', 'Croissant - Spinach Feta', '5', '59.91', '299.55'),
    @('9090', 'Sugared Egg Yolk', '1', '119.40', '119.40')
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $rng = $ws.Range("A" + $row + ":E" + $row)
    $rng.NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}
